$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the current row 353, shifting the existing
# rows 353:367 down to 355:369 (dimension grows from T367 to T369).
$ws.Rows("353:354").Insert()

# ---- New row 353 ----
$ws.Cells.Item(353,1).Value  = 10
$ws.Cells.Item(353,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(353,3).Value  = "La Araucanía"
$ws.Cells.Item(353,4).Value  = 44509
$ws.Cells.Item(353,5).Value  = 9
$ws.Cells.Item(353,6).Value  = "Fruta"
$ws.Cells.Item(353,7).Value  = 100108
$ws.Cells.Item(353,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(353,9).Value  = 100108006
$ws.Cells.Item(353,10).Value = "Plátano"
$ws.Cells.Item(353,11).Value = "Barraganete"
$ws.Cells.Item(353,12).Value = "Primera"
$ws.Cells.Item(353,13).Value = 55
$ws.Cells.Item(353,14).Value = 30000
$ws.Cells.Item(353,15).Value = 30000
$ws.Cells.Item(353,16).Value = 30000
$ws.Cells.Item(353,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(353,18).Value = "Ecuador"
$ws.Cells.Item(353,19).Value = 1500
$ws.Cells.Item(353,20).Value = 20

# ---- New row 354 ----
$ws.Cells.Item(354,1).Value  = 10
$ws.Cells.Item(354,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(354,3).Value  = "La Araucanía"
$ws.Cells.Item(354,4).Value  = 44509
$ws.Cells.Item(354,5).Value  = 9
$ws.Cells.Item(354,6).Value  = "Fruta"
$ws.Cells.Item(354,7).Value  = 100108
$ws.Cells.Item(354,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(354,9).Value  = 100108006
$ws.Cells.Item(354,10).Value = "Plátano"
$ws.Cells.Item(354,11).Value = "Sin especificar"
$ws.Cells.Item(354,12).Value = "Pintón"
$ws.Cells.Item(354,13).Value = 450
$ws.Cells.Item(354,14).Value = 19000
$ws.Cells.Item(354,15).Value = 21000
$ws.Cells.Item(354,16).Value = 20111
$ws.Cells.Item(354,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(354,18).Value = "Ecuador"
$ws.Cells.Item(354,19).Value = 1006
$ws.Cells.Item(354,20).Value = 20
